$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Purchase 22-23: remove the "Shruti Infotech" line item (old Sr.No 5, rows 13) ---
$ws1.Rows.Item(13).Delete()
$ws1.Rows.Item(13).Delete()

# --- Insert a new blank row before the (shifted) blank separator row 15 so the
#     "Shree Laxmi Lighting Hub" group gets a third invoice line ---
$ws1.Rows.Item(15).Insert()

# Pick up matching number/border formatting from the row above (same group)
$ws1.Range("A13:F13").Copy()
$ws1.Range("A15:F15").PasteSpecial(-4122)
$ws1.Rows.Item(15).RowHeight = 14.4

# New third "Shree Laxmi Lighting Hub" invoice (SLH/3668)
$ws1.Range("B15").Value = 45311
$ws1.Range("C15").Value = "SLH/3668"
$ws1.Range("D15").Value = "Shree Laxmi Lighting Hub"
$ws1.Range("E15").Value = 1756

# The running-total formula now lives on the new last row of the group
$ws1.Range("F14").ClearContents()
$ws1.Range("F15").Formula = "=E13+E14+E15"

# Renumber the Sr. No column for the groups that shifted up
$ws1.Range("A13").Value = 5
$ws1.Range("A17").Value = 6

# --- Selection / active sheet bookkeeping: Purchase 22-23 becomes the active tab ---
$ws1.Activate()
$ws1.Range("A19").Select()
